# Update the work shift for "Long" on Monday, Tuesday and Wednesday
# (columns B, C, D) from "08:00 - 17:00" to "08:00 - 16:00" for every
# week block in the schedule (rows 3, 9, 15, 21, 27, 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 9, 15, 21, 27, 33)
$cols = @("B", "C", "D")

foreach ($row in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = "08:00 - 16:00"
    }
}
